$wb = $excel.ActiveWorkbook

$wsStudies = $wb.Worksheets.Item("studies")
$wsStudies.Range("A1").Value = "study_id"

$wsSurveys = $wb.Worksheets.Item("surveys")
$wsSurveys.Range("B1").Value = "survey_id"
$wsSurveys.Range("E1").Value = "latitude"
$wsSurveys.Range("F1").Value = "longitude"

# remember the currently active sheet so we can restore it
$prevActive = $wb.ActiveSheet.Name

$wsStudies.Range("A2").Select()

# restore original active sheet
$wb.Worksheets.Item($prevActive).Activate()
